$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to text format first so values like
# "1.00" or "6.60" are not silently coerced into numbers by Excel,
# matching the inline-string (text) storage used in the workbook.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.957.83'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '1.558.42'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '207.99'
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("D6").Value = '0.489'
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '22.11'
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = '0.248'
$ws.Range("E9").Value = '  +0.49%  '
$ws.Range("D10").Value = '0.0598'
$ws.Range("E10").Value = '  +1.66%  '
$ws.Range("D11").Value = '0.0855'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '1.781.01'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '1.512.76'
$ws.Range("E13").Value = '  -2.79%  '
$ws.Range("D14").Value = '3.75'
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("D15").Value = '0.520'
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").Value = '61.91'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").Value = '26.956.91'
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").Value = '0.0₃0707'
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("D19").Value = '215.97'
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("E20").Value = '  +0.67%  '
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = '4.11'
$ws.Range("E22").Value = '  +1.34%  '
$ws.Range("D23").Value = '9.22'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("E24").Value = '  -1.09%  '
$ws.Range("D25").Value = '153.02'
$ws.Range("E25").Value = '  -0.72%  '
$ws.Range("D26").Value = '6.60'
$ws.Range("E26").Value = '  -0.93%  '
$ws.Range("D27").Value = '15.08'
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("E28").Value = '  +1.03%  '
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("E31").Value = '  +3.37%  '
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("D33").Value = '3.16'
$ws.Range("E33").Value = '  +2.94%  '
$ws.Range("D34").Value = '1.423.18'
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("D36").Value = '1.06'
$ws.Range("E36").Value = '  +8.64%  '
$ws.Range("E37").Value = '  +1.97%  '
$ws.Range("D38").Value = '0.0166'
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("D39").Value = '0.532'
$ws.Range("E39").Value = '  +2.30%  '
$ws.Range("E40").Value = '  +2.60%  '
$ws.Range("D41").Value = '0.807'
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +2.20%  '
$ws.Range("D44").Value = '2.31'
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("D45").Value = '64.55'
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").Value = '1.694.05'
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("D48").Value = '87.16'
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0102'
$ws.Range("E49").Value = '  +2.33%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.0519'
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").Value = '0.0959'
$ws.Range("E51").Value = '  -0.24%  '

# Restore default (General/Normal) styling so no style index is left on
# these cells - matches the original workbook which had no "s" attribute
# on these cells.
$ws.Range("B2:E51").Style = "Normal"
